# Working hours workbook update:
# - Insert a new data row at row 54 (pushes the summary rows down by one).
# - Correct the end time recorded in (the now shifted) row 53.
# - Fill the newly inserted row 54 with the new time entry and its formulas.
# - Re-point the active selection to A55 (first cell of the now-empty
#   placeholder row that used to sit at row 54).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the old row 54 (the empty placeholder row).
# Excel automatically shifts all rows below down by one and adjusts the
# formulas/shared-formula ranges that reference them.
$ws.Range("A54").EntireRow.Insert()

# The end time of the entry that is now on row 53 changes from 18:30 to 18:45
$ws.Range("E53").Value = 0.78125

# Populate the newly inserted row 54 with the new working-hours entry.
$ws.Range("A54").Value = 2014
$ws.Range("B54").Value = 5
$ws.Range("C54").Value = 6
$ws.Range("D54").Value = 0.83333333333333337
$ws.Range("E54").Value = 0.91666666666666663
$ws.Range("F54").Formula = "=(E54-D54)*24*60"
$ws.Range("G54").Formula = "=F54/60"

# Match the workbook's saved selection state.
$ws.Range("A55").Select()
